# Insert a new data row at row 160 (shifts existing rows 160-222 down to 161-223)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(160).Insert()

$ws.Cells.Item(160, 1).Value = 3
$ws.Cells.Item(160, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44875
$ws.Cells.Item(160, 5).Value = 5
$ws.Cells.Item(160, 6).Value = 100112026
$ws.Cells.Item(160, 7).Value = "Haba"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 90
$ws.Cells.Item(160, 11).Value = 8000
$ws.Cells.Item(160, 12).Value = 8500
$ws.Cells.Item(160, 13).Value = 8278
$ws.Cells.Item(160, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(160, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(160, 16).Value = 331
$ws.Cells.Item(160, 17).Value = 25
$ws.Cells.Item(160, 18).Value = "Hortaliza"
